# Append a new log row to the "統計" (statistics) sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("統計")

# Find the next empty row after the existing data (row 11 -> 12)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = "2025-08-28T12:44:47.229040"
$ws.Cells.Item($newRow, 2).Value = 13
$ws.Cells.Item($newRow, 3).Value = "全案件リスト"
$ws.Cells.Item($newRow, 4).Value = 53.8
$ws.Cells.Item($newRow, 5).Value = 6
$ws.Cells.Item($newRow, 6).Value = 4
$ws.Cells.Item($newRow, 7).Value = 13
